# Fix: prevent hidden columns from being labeled upon detecting changes.
# Column K/V ("Bedingung_FV2410"/"Bedingung_FV2504") is a hidden column that
# should be skipped when diffing rows. Rows whose ONLY detected difference
# was that hidden column were incorrectly flagged "AENDERUNG" in column L.
# This clears those false-positive flags, and restores the "new group"
# highlight styling on rows whose flag removal leaves them with no other
# visual marker of being the first row of a new lfd. Position group.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows where only the L (Aenderung) cell needs its false-positive flag removed.
$rowsLOnly = @(76,77,78,80,81,82,84,85,87,88,89,91,92,93,94,96,97,98,100,101,102,104,105,107,108)

foreach ($r in $rowsLOnly) {
    $ws.Range("L2").Copy()
    $ws.Range("L$r").PasteSpecial(-4122)
    $ws.Range("L$r").ClearContents()
}

# Rows that start a new "lfd. Position" group whose only change was the
# hidden-column false positive: once the flag is removed they need the
# group-header formatting (style used by row 2, a known-good group header)
# re-applied across the whole row so the new-group boundary stays visible.
$rowsFullRewrite = @(83,86,90,95,99,103,106)

foreach ($r in $rowsFullRewrite) {
    $ws.Range("A2:V2").Copy()
    $ws.Range("A$r`:V$r").PasteSpecial(-4122)
    $ws.Range("L$r").ClearContents()
}

$excel.CutCopyMode = 0
